# aggiornamento fino a 6 gennaio 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(44539,44540,44541,44542,44543,44544,44545,44546,44547,44548,
           44550,44551,44552,44553,44554,44555,44556,44557,44558,44559,
           44560,44561,44562,44563,44564,44565,44566)
$colB  = @(6,1,0,12,4,0,0,3,0,1,0,1,0,0,3,0,1,0,0,1,0,0,1,1,0,0,1)
$colC  = @(14,12,12,21,24,24,23,20,19,20,8,5,5,5,5,5,5,5,4,5,5,2,3,3,3,3,3)
$colD  = @(1167.639699749791,1000.834028356964,1000.834028356964,1751.459549624687,
           2001.668056713928,2001.668056713928,1918.265221017515,1668.056713928273,
           1584.65387823186,1668.056713928273,667.2226855713094,417.0141784820684,
           417.0141784820684,417.0141784820684,417.0141784820684,417.0141784820684,
           417.0141784820684,417.0141784820684,333.6113427856547,417.0141784820684,
           417.0141784820684,166.8056713928273,250.208507089241,250.208507089241,
           250.208507089241,250.208507089241,250.208507089241)

$startRow = 465
$endRow = $startRow + $dates.Length - 1

# Clone the date column's cell formatting (style index "2" in the source file)
# from the last existing data row down onto the new rows.
$ws.Range("A464").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}
